# Weekly update of the "Hortaliza, Vega Monumental Concepción - Ají" sheet.
# A new weekly record is inserted at row 147 (pushing the existing rows
# 147:158 down to 148:159), and the new row is filled with the latest
# observation for the market.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 147; this shifts the existing
# rows 147-158 down to 148-159, carrying their formatting (including the
# date style on column D) along with them.
$ws.Rows.Item(147).Insert()

# Populate the newly inserted row 147 with the new weekly record.
$ws.Range("A147").Value = 11
$ws.Range("B147").Value = "Vega Monumental Concepción"
$ws.Range("C147").Value = "Bíobío"
$ws.Range("D147").Value = 44931
$ws.Range("E147").Value = 8
$ws.Range("F147").Value = 100112021
$ws.Range("G147").Value = "Ají"
$ws.Range("H147").Value = "Americana (o)"
$ws.Range("I147").Value = "Primera"
$ws.Range("J147").Value = 270
$ws.Range("K147").Value = 20000
$ws.Range("L147").Value = 22000
$ws.Range("M147").Value = 21111
$ws.Range("N147").Value = "$/caja 15 kilos"
$ws.Range("O147").Value = "Región Metropolitana"
$ws.Range("P147").Value = 1407
$ws.Range("Q147").Value = 15
$ws.Range("R147").Value = "Hortaliza"
